$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 43-46 as "Purchased" (column N) using the built-in "Good" cell style,
# matching the highlighted look used elsewhere for the Purchased? column.
# Build the combined format (Good style + left/center alignment) on a scratch
# cell first, then paste just the formatting onto the target cells - this
# avoids leaving behind unused intermediate cell styles.
$fmtSrc = $ws.Range("ZZ1")
$fmtSrc.Style = "Good"
$fmtSrc.HorizontalAlignment = -4131
$fmtSrc.VerticalAlignment = -4108

$markRange = $ws.Range("N43:N46")
$markRange.Value = "y"
$fmtSrc.Copy()
$markRange.PasteSpecial(-4122)
$fmtSrc.Clear()
$excel.CutCopyMode = 0

# Add the new "Thonk" order line in row 59 (previously a blank gap row),
# copying the number format from the row above for the price cell.
$ws.Range("A59").Value = "Thonk"
$ws.Range("B59").Value = 22.74
$ws.Range("B58").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Remove the old blank gap row (was row 59) that is now row 60, and move the
# "Total" row down from 60 to 61, updating the sum range to include B59.
$ws.Range("A60:B60").Clear()
$ws.Range("A61").Value = "Total"
$ws.Range("B61").Formula = "=SUM(B56:B59)"

# Restore the selection/scroll state roughly to what was saved in the workbook.
$ws.Activate()
$ws.Range("P32").Select()
